# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (cols H-N) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1177.7778
$ws.Range("I70").Value = 1166.6666
$ws.Range("K70").Value = 3499.9998
$ws.Range("M70").Value = -3229.9998

$ws.Range("H73").Value = 1177.7778
$ws.Range("I73").Value = 1166.6666
$ws.Range("K73").Value = 3499.9998
$ws.Range("M73").Value = -2563.9998

$ws.Range("H129").Value = 1020
$ws.Range("J129").Value = 1178.6097
$ws.Range("L129").Value = 3535.8291
$ws.Range("N129").Value = -13535.8291

$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040

$ws.Range("H137").Value = 4584
$ws.Range("I137").Value = 4584
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 13752
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -11202
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 200167.4
$ws.Range("I138").Value = 2742.6667
$ws.Range("J138").Value = 274201.7
$ws.Range("K138").Value = 8228.000100000001
$ws.Range("L138").Value = 822605.1000000001
$ws.Range("M138").Value = -3088.000100000001
$ws.Range("N138").Value = -832885.1000000001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 13512.75
$ws.Range("J23").Value = 9746.532999999999
$ws.Range("L23").Value = 9746.532999999999
$ws.Range("N23").Value = -10264.533

$ws.Range("H32").Value = 444171.16
$ws.Range("I32").Value = 506164.9
$ws.Range("K32").Value = 506164.9
$ws.Range("M32").Value = -505877.9

$ws.Range("H37").Value = 11896.435
$ws.Range("J37").Value = 11896.435
$ws.Range("L37").Value = 11896.435
$ws.Range("N37").Value = -12442.435

$ws.Range("H44").Value = 29332
$ws.Range("J44").Value = 29332
$ws.Range("L44").Value = 29332
$ws.Range("N44").Value = -30308

$ws.Range("H63").Value = 5104.5625
$ws.Range("I63").Value = 3635.75
$ws.Range("J63").Value = 6573.375
$ws.Range("K63").Value = 3635.75
$ws.Range("L63").Value = 6573.375
$ws.Range("M63").Value = -2949.75
$ws.Range("N63").Value = -7945.375

$ws.Range("H66").Value = 5104.5625
$ws.Range("I66").Value = 3635.75
$ws.Range("J66").Value = 6573.375
$ws.Range("K66").Value = 18178.75
$ws.Range("L66").Value = 32866.875
$ws.Range("M66").Value = -14746.75
$ws.Range("N66").Value = -39730.875

$ws.Range("H132").Value = 3229.8057
$ws.Range("I132").Value = 3077.75
$ws.Range("K132").Value = 9233.25
$ws.Range("M132").Value = -6703.25


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 24999
$ws.Range("J35").Value = 24999
$ws.Range("L35").Value = 24999
$ws.Range("N35").Value = -25619

$ws.Range("H82").Value = 15445.833
$ws.Range("J82").Value = 30099.25
$ws.Range("L82").Value = 30099.25
$ws.Range("N82").Value = -30865.25

$ws.Range("H85").Value = 15445.833
$ws.Range("J85").Value = 30099.25
$ws.Range("L85").Value = 30099.25
$ws.Range("N85").Value = -32751.25


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 92.14286
$ws.Range("I7").Value = 49
$ws.Range("K7").Value = 49
$ws.Range("M7").Value = 64

$ws.Range("H31").Value = 2438.923
$ws.Range("I31").Value = 975.9666999999999
$ws.Range("J31").Value = 7315.4443
$ws.Range("K31").Value = 975.9666999999999
$ws.Range("L31").Value = 7315.4443
$ws.Range("M31").Value = -680.9666999999999
$ws.Range("N31").Value = -7905.4443

$ws.Range("H34").Value = 2438.923
$ws.Range("I34").Value = 975.9666999999999
$ws.Range("J34").Value = 7315.4443
$ws.Range("K34").Value = 975.9666999999999
$ws.Range("L34").Value = 7315.4443
$ws.Range("M34").Value = -773.9666999999999
$ws.Range("N34").Value = -7719.4443

$ws.Range("H41").Value = 17713.285
$ws.Range("J41").Value = 17713.285
$ws.Range("L41").Value = 17713.285
$ws.Range("N41").Value = -18569.285

$ws.Range("H50").Value = 19249
$ws.Range("J50").Value = 19249
$ws.Range("L50").Value = 19249
$ws.Range("N50").Value = -20499

$ws.Range("H51").Value = 17766
$ws.Range("J51").Value = 17766
$ws.Range("L51").Value = 17766
$ws.Range("N51").Value = -19238

$ws.Range("H59").Value = 28199.2
$ws.Range("J59").Value = 28199.2
$ws.Range("L59").Value = 28199.2
$ws.Range("N59").Value = -30489.2

$ws.Range("H60").Value = 10376.125
$ws.Range("J60").Value = 10376.125
$ws.Range("L60").Value = 10376.125
$ws.Range("N60").Value = -11398.125

$ws.Range("H61").Value = 17766
$ws.Range("J61").Value = 17766
$ws.Range("L61").Value = 17766
$ws.Range("N61").Value = -18462

$ws.Range("H68").Value = 27582.666
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 27582.666
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 27582.666
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -29080.666

$ws.Range("H71").Value = 27582.666
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 27582.666
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 82747.99800000001
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -90235.99800000001

$ws.Range("H74").Value = 25942.5
$ws.Range("J74").Value = 27391.818
$ws.Range("L74").Value = 27391.818
$ws.Range("N74").Value = -29139.818

$ws.Range("H77").Value = 25942.5
$ws.Range("J77").Value = 27391.818
$ws.Range("L77").Value = 82175.454
$ws.Range("N77").Value = -90911.454

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 4000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 12000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -12338

$ws.Range("H45").Value = 1051.2
$ws.Range("J45").Value = 1051.2
$ws.Range("L45").Value = 3153.6
$ws.Range("N45").Value = -4217.6

$ws.Range("H55").Value = 2300
$ws.Range("J55").Value = 2300
$ws.Range("L55").Value = 6900
$ws.Range("N55").Value = -7254

$ws.Range("H60").Value = 2774.8333
$ws.Range("I60").Value = 490
$ws.Range("J60").Value = 2927.1555
$ws.Range("K60").Value = 1470
$ws.Range("L60").Value = 8781.466499999999
$ws.Range("M60").Value = -1219
$ws.Range("N60").Value = -9283.466499999999

$ws.Range("H70").Value = 2825
$ws.Range("I70").Value = 650
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 1950
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -1635
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 2825
$ws.Range("I73").Value = 650
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 1950
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -858
$ws.Range("N73").Value = -17184

$ws.Range("H81").Value = 4124.75
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 5166.3335
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 15499.0005
$ws.Range("M81").Value = -1877
$ws.Range("N81").Value = -17745.0005

$ws.Range("H84").Value = 4124.75
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 5166.3335
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 46497.0015
$ws.Range("M84").Value = -3384
$ws.Range("N84").Value = -57729.0015

$ws.Range("H123").Value = 7785.7144
$ws.Range("I123").Value = 2750
$ws.Range("J123").Value = 9800
$ws.Range("K123").Value = 8250
$ws.Range("L123").Value = 29400
$ws.Range("M123").Value = -5800
$ws.Range("N123").Value = -34300

$ws.Range("H124").Value = 1881.5454
$ws.Range("I124").Value = 1495
$ws.Range("J124").Value = 1967.4445
$ws.Range("K124").Value = 4485
$ws.Range("L124").Value = 5902.333500000001
$ws.Range("M124").Value = 425
$ws.Range("N124").Value = -15722.3335

$ws.Range("H129").Value = 2486.6667
$ws.Range("I129").Value = 1000
$ws.Range("J129").Value = 2621.818
$ws.Range("K129").Value = 3000
$ws.Range("L129").Value = 7865.454000000001
$ws.Range("M129").Value = 2000
$ws.Range("N129").Value = -17865.454


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 7642375.5
$ws.Range("I10").Value = 12001100
$ws.Range("J10").Value = 377834.66
$ws.Range("K10").Value = 12001100
$ws.Range("L10").Value = 377834.66
$ws.Range("M10").Value = -12000931
$ws.Range("N10").Value = -378172.66

$ws.Range("H100").Value = 40000
$ws.Range("J100").Value = 40000
$ws.Range("L100").Value = 40000
$ws.Range("N100").Value = -42164

$ws.Range("H130").Value = 55400
$ws.Range("J130").Value = 55400
$ws.Range("L130").Value = 55400
$ws.Range("N130").Value = -65440

$ws.Range("H131").Value = 34500
$ws.Range("J131").Value = 34500
$ws.Range("L131").Value = 34500
$ws.Range("N131").Value = -44580


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 23519.143
$ws.Range("J2").Value = 39751
$ws.Range("L2").Value = 39751
$ws.Range("N2").Value = -39975

$ws.Range("H100").Value = 3098.3333
$ws.Range("I100").Value = 2900
$ws.Range("J100").Value = 3495
$ws.Range("K100").Value = 2900
$ws.Range("L100").Value = 3495
$ws.Range("M100").Value = -2359
$ws.Range("N100").Value = -4577


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 15867.333
$ws.Range("J101").Value = 15867.333
$ws.Range("L101").Value = 15867.333
$ws.Range("N101").Value = -22357.333

